$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed ligand/receptor/edge expression-specificity values using the updated TPM matrix.
# Columns: G=Ligand avg expr, H=Ligand total expr, I/J=Ligand specificity (avg/total),
#          M=Receptor avg expr, N=Receptor total expr, O/P=Receptor specificity (avg/total),
#          Q=Edge avg weight, R=Edge total weight, S/T=Edge specificity (avg/total)

# Row 2
$ws.Cells.Item(2, 7).Value = 1.576331333333333
$ws.Cells.Item(2, 8).Value = 4.728994
$ws.Cells.Item(2, 9).Value = 0.01463337290888519
$ws.Cells.Item(2, 10).Value = 0.01463337290888519
$ws.Cells.Item(2, 13).Value = 61.04160633333334
$ws.Cells.Item(2, 14).Value = 183.124819
$ws.Cells.Item(2, 15).Value = 0.2043613460574534
$ws.Cells.Item(2, 16).Value = 0.2043613460574534
$ws.Cells.Item(2, 17).Value = 96.22179670023179
$ws.Cells.Item(2, 18).Value = 865.9961703020862
$ws.Cells.Item(2, 19).Value = 0.00299049578502045
$ws.Cells.Item(2, 20).Value = 0.00299049578502045

# Row 3
$ws.Cells.Item(3, 7).Value = 1.576331333333333
$ws.Cells.Item(3, 8).Value = 4.728994
$ws.Cells.Item(3, 9).Value = 0.01463337290888519
$ws.Cells.Item(3, 10).Value = 0.01463337290888519
$ws.Cells.Item(3, 15).Value = 0.3559304658284363
$ws.Cells.Item(3, 16).Value = 0.3559304658284363
$ws.Cells.Item(3, 17).Value = 167.5868239424013
$ws.Cells.Item(3, 18).Value = 1508.281415481612
$ws.Cells.Item(3, 19).Value = 0.005208463236100725
$ws.Cells.Item(3, 20).Value = 0.005208463236100727

# Row 4
$ws.Cells.Item(4, 7).Value = 1.576331333333333
$ws.Cells.Item(4, 8).Value = 4.728994
$ws.Cells.Item(4, 9).Value = 0.01463337290888519
$ws.Cells.Item(4, 10).Value = 0.01463337290888519
$ws.Cells.Item(4, 13).Value = 131.3384093333333
$ws.Cells.Item(4, 14).Value = 394.015228
$ws.Cells.Item(4, 15).Value = 0.4397081881141102
$ws.Cells.Item(4, 16).Value = 0.4397081881141103
$ws.Cells.Item(4, 17).Value = 207.0328499022924
$ws.Cells.Item(4, 18).Value = 1863.295649120632
$ws.Cells.Item(4, 19).Value = 0.006434413887764012
$ws.Cells.Item(4, 20).Value = 0.006434413887764014

# Row 5
$ws.Cells.Item(5, 7).Value = 74.31489566666666
$ws.Cells.Item(5, 9).Value = 0.6898788078237544
$ws.Cells.Item(5, 10).Value = 0.6898788078237544
$ws.Cells.Item(5, 13).Value = 61.04160633333334
$ws.Cells.Item(5, 14).Value = 183.124819
$ws.Cells.Item(5, 15).Value = 0.2043613460574534
$ws.Cells.Item(5, 16).Value = 0.2043613460574534
$ws.Cells.Item(5, 17).Value = 4536.300605987406
$ws.Cells.Item(5, 18).Value = 40826.70545388666
$ws.Cells.Item(5, 19).Value = 0.1409845617833737
$ws.Cells.Item(5, 20).Value = 0.1409845617833737

# Row 6
$ws.Cells.Item(6, 7).Value = 74.31489566666666
$ws.Cells.Item(6, 9).Value = 0.6898788078237544
$ws.Cells.Item(6, 10).Value = 0.6898788078237544
$ws.Cells.Item(6, 15).Value = 0.3559304658284363
$ws.Cells.Item(6, 16).Value = 0.3559304658284363
$ws.Cells.Item(6, 17).Value = 7900.74844864738
$ws.Cells.Item(6, 19).Value = 0.2455488854338752
$ws.Cells.Item(6, 20).Value = 0.2455488854338752

# Row 7
$ws.Cells.Item(7, 7).Value = 74.31489566666666
$ws.Cells.Item(7, 9).Value = 0.6898788078237544
$ws.Cells.Item(7, 10).Value = 0.6898788078237544
$ws.Cells.Item(7, 13).Value = 131.3384093333333
$ws.Cells.Item(7, 14).Value = 394.015228
$ws.Cells.Item(7, 15).Value = 0.4397081881141102
$ws.Cells.Item(7, 16).Value = 0.4397081881141103
$ws.Cells.Item(7, 17).Value = 9760.400186632623
$ws.Cells.Item(7, 18).Value = 87843.60167969363
$ws.Cells.Item(7, 19).Value = 0.3033453606065055
$ws.Cells.Item(7, 20).Value = 0.3033453606065055

# Row 8
$ws.Cells.Item(8, 7).Value = 31.83044066666666
$ws.Cells.Item(8, 8).Value = 95.491322
$ws.Cells.Item(8, 9).Value = 0.2954878192673605
$ws.Cells.Item(8, 10).Value = 0.2954878192673605
$ws.Cells.Item(8, 13).Value = 61.04160633333334
$ws.Cells.Item(8, 14).Value = 183.124819
$ws.Cells.Item(8, 15).Value = 0.2043613460574534
$ws.Cells.Item(8, 16).Value = 0.2043613460574534
$ws.Cells.Item(8, 17).Value = 1942.981228591191
$ws.Cells.Item(8, 18).Value = 17486.83105732072
$ws.Cells.Item(8, 19).Value = 0.06038628848905932
$ws.Cells.Item(8, 20).Value = 0.06038628848905932

# Row 9
$ws.Cells.Item(9, 7).Value = 31.83044066666666
$ws.Cells.Item(9, 8).Value = 95.491322
$ws.Cells.Item(9, 9).Value = 0.2954878192673605
$ws.Cells.Item(9, 10).Value = 0.2954878192673605
$ws.Cells.Item(9, 15).Value = 0.3559304658284363
$ws.Cells.Item(9, 16).Value = 0.3559304658284363
$ws.Cells.Item(9, 17).Value = 3384.03630202135
$ws.Cells.Item(9, 18).Value = 30456.32671819215
$ws.Cells.Item(9, 19).Value = 0.1051731171584604
$ws.Cells.Item(9, 20).Value = 0.1051731171584604

# Row 10
$ws.Cells.Item(10, 7).Value = 31.83044066666666
$ws.Cells.Item(10, 8).Value = 95.491322
$ws.Cells.Item(10, 9).Value = 0.2954878192673605
$ws.Cells.Item(10, 10).Value = 0.2954878192673605
$ws.Cells.Item(10, 13).Value = 131.3384093333333
$ws.Cells.Item(10, 14).Value = 394.015228
$ws.Cells.Item(10, 15).Value = 0.4397081881141102
$ws.Cells.Item(10, 16).Value = 0.4397081881141103
$ws.Cells.Item(10, 17).Value = 4180.559445539046
$ws.Cells.Item(10, 18).Value = 37625.03500985141
$ws.Cells.Item(10, 19).Value = 0.1299284136198408
$ws.Cells.Item(10, 20).Value = 0.1299284136198408
